# Updates cryptos list cell values per the Dec 23 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    # Leading apostrophe forces Excel to store the value as text even when it
    # looks numeric (e.g. "96.52"); Style is reset afterwards so no numeric/text
    # format gets stamped onto the cell (matches original unstyled data cells).
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "43.747.25"
Set-TextCell $ws.Range("E2") "  +0.43%  "
Set-TextCell $ws.Range("D3") "2.292.65"
Set-TextCell $ws.Range("E3") "  -0.65%  "
Set-TextCell $ws.Range("E4") "  -0.21%  "
Set-TextCell $ws.Range("D5") "96.52"
Set-TextCell $ws.Range("E5") "  +3.62%  "
Set-TextCell $ws.Range("D6") "269.95"
Set-TextCell $ws.Range("E6") "  +1.05%  "
Set-TextCell $ws.Range("D7") "0.617"
Set-TextCell $ws.Range("E7") "  -1.14%  "
Set-TextCell $ws.Range("E8") "  -0.17%  "
Set-TextCell $ws.Range("D9") "0.608"
Set-TextCell $ws.Range("E9") "  -0.95%  "
Set-TextCell $ws.Range("E10") "  +3.05%  "
Set-TextCell $ws.Range("D11") "0.0935"
Set-TextCell $ws.Range("E11") "  -0.18%  "
Set-TextCell $ws.Range("E12") "  -1.43%  "
Set-TextCell $ws.Range("E13") "  +1.79%  "
Set-TextCell $ws.Range("D14") "15.76"
Set-TextCell $ws.Range("E14") "  +2.83%  "
Set-TextCell $ws.Range("D15") "2.636.99"
Set-TextCell $ws.Range("E15") "  -0.75%  "
Set-TextCell $ws.Range("D16") "0.855"
Set-TextCell $ws.Range("E16") "  +0.34%  "
Set-TextCell $ws.Range("D17") "2.297.85"
Set-TextCell $ws.Range("E17") "  -0.74%  "
Set-TextCell $ws.Range("D18") "43.692.54"
Set-TextCell $ws.Range("E18") "  +0.37%  "
Set-TextCell $ws.Range("D19") "0.0000111"
Set-TextCell $ws.Range("E19") "  +4.21%  "
Set-TextCell $ws.Range("D20") "6.19"
Set-TextCell $ws.Range("E20") "  -2.11%  "
Set-TextCell $ws.Range("D21") "72.16"
Set-TextCell $ws.Range("E21") "  +1.29%  "
Set-TextCell $ws.Range("D22") "2.51"
Set-TextCell $ws.Range("E22") "  +11.40%  "
Set-TextCell $ws.Range("D23") "232.71"
Set-TextCell $ws.Range("E23") "  -1.64%  "
Set-TextCell $ws.Range("D24") "9.15"
Set-TextCell $ws.Range("E24") "  -4.61%  "
Set-TextCell $ws.Range("D25") "2.74"
Set-TextCell $ws.Range("E25") "  +10.09%  "
Set-TextCell $ws.Range("E26") "  -0.10%  "
Set-TextCell $ws.Range("D27") "11.31"
Set-TextCell $ws.Range("E27") "  +0.59%  "
Set-TextCell $ws.Range("E28") "  -1.97%  "
Set-TextCell $ws.Range("D29") "38.71"
Set-TextCell $ws.Range("E29") "  +1.04%  "
Set-TextCell $ws.Range("E30") "  -2.61%  "
Set-TextCell $ws.Range("B31") "EthereumClassic"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D31") "22.27"
Set-TextCell $ws.Range("E31") "  -0.94%  "
Set-TextCell $ws.Range("B32") "Monero"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D32") "174.99"
Set-TextCell $ws.Range("E32") "  +2.08%  "
Set-TextCell $ws.Range("D33") "0.0899"
Set-TextCell $ws.Range("E33") "  +1.07%  "
Set-TextCell $ws.Range("E34") "  -0.92%  "
Set-TextCell $ws.Range("E35") "  +0.61%  "
Set-TextCell $ws.Range("D36") "4.51"
Set-TextCell $ws.Range("E36") "  +4.06%  "
Set-TextCell $ws.Range("E37") "  -0.30%  "
Set-TextCell $ws.Range("D38") "0.106"
Set-TextCell $ws.Range("E38") "  -2.04%  "
Set-TextCell $ws.Range("D39") "3.49"
Set-TextCell $ws.Range("E39") "  +3.16%  "
Set-TextCell $ws.Range("D40") "0.236"
Set-TextCell $ws.Range("E40") "  +1.73%  "
Set-TextCell $ws.Range("E41") "  +0.20%  "
Set-TextCell $ws.Range("D42") "12.27"
Set-TextCell $ws.Range("E42") "  +2.47%  "
Set-TextCell $ws.Range("E43") "  -0.57%  "
Set-TextCell $ws.Range("D44") "64.28"
Set-TextCell $ws.Range("E44") "  +4.62%  "
Set-TextCell $ws.Range("E45") "  -2.83%  "
Set-TextCell $ws.Range("E46") "  -3.77%  "
Set-TextCell $ws.Range("E47") "  +0.21%  "
Set-TextCell $ws.Range("E48") "  -0.47%  "
Set-TextCell $ws.Range("E49") "  -2.15%  "
Set-TextCell $ws.Range("D50") "1.56"
Set-TextCell $ws.Range("E50") "  +15.42%  "
Set-TextCell $ws.Range("D51") "0.434"
Set-TextCell $ws.Range("E51") "  +4.50%  "
